$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.750.91"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "1.660.37"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").Value = "'302.96"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "'0.3824"
$ws.Range("D8").Value = "'0.3619"
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").Value = "'51.26"
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("D10").Value = "'0.08196"
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").Value = "'1.232"
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").Value = "'22.58"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").Value = "'6.472"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").Value = "'7.432"
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("E16").Value = "  -1.13%  "
$ws.Range("D17").Value = "1.660.44"
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("D18").Value = "'97.85"
$ws.Range("E18").Value = "  +2.71%  "
$ws.Range("D19").Value = "'0.07021"
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("D20").Value = "'6.823"
$ws.Range("E20").Value = "  +3.81%  "
$ws.Range("D21").Value = "'17.65"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "'12.84"
$ws.Range("E23").Value = "  +2.44%  "
$ws.Range("D24").Value = "23.752.01"
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("D25").Value = "'2.506"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("D26").Value = "'3.002"
$ws.Range("E26").Value = "  -2.59%  "
$ws.Range("D27").Value = "'21.22"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").Value = "'152.79"
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("D29").Value = "'5.232"
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("D30").Value = "'134.26"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("D31").Value = "1.843.91"
$ws.Range("E31").Value = "  +1.45%  "
$ws.Range("D32").Value = "'7.139"
$ws.Range("E32").Value = "  +7.77%  "
$ws.Range("D33").Value = "'2.240"
$ws.Range("E33").Value = "  +4.31%  "
$ws.Range("D34").Value = "'12.05"
$ws.Range("E34").Value = "  +4.42%  "
$ws.Range("D35").Value = "'1.056"
$ws.Range("E35").Value = "  -3.46%  "
$ws.Range("D36").Value = "'0.02817"
$ws.Range("E36").Value = "  +1.86%  "
$ws.Range("D37").Value = "'0.2518"
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("D38").Value = "'0.08814"
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("D39").Value = "'6.091"
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("D40").Value = "'0.07018"
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("D41").Value = "'13.03"
$ws.Range("E41").Value = "  +5.63%  "
$ws.Range("D42").Value = "'0.7009"
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("D43").Value = "'1.335"
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("D44").Value = "'16.02"
$ws.Range("E44").Value = "  +2.61%  "
$ws.Range("D45").Value = "'0.6523"
$ws.Range("E45").Value = "  -0.65%  "
$ws.Range("D46").Value = "'0.9997"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").Value = "'2.308"
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("D48").Value = "'3.967"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("D50").Value = "'128.50"
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("D51").Value = "'1.183"
$ws.Range("E51").Value = "  -1.27%  "
